$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("H4").Value = 3.8

# Row 6
$ws.Range("H6").Value = 4.65
$ws.Range("I6").Value = 4.6
$ws.Range("L6").Value = 1.08
$ws.Range("M6").Value = 6.6
$ws.Range("N6").Value = 1.26
$ws.Range("O6").Value = 3.5
$ws.Range("P6").Value = 1.17
$ws.Range("Q6").Value = 4.5
$ws.Range("X6").Value = 10.75
$ws.Range("Y6").Value = 14.5
$ws.Range("AB6").Value = 12.5
$ws.Range("AC6").Value = 28
$ws.Range("AD6").Value = 110
$ws.Range("AE6").Value = 29
$ws.Range("AF6").Value = 40
$ws.Range("AH6").Value = 80
$ws.Range("AJ6").Value = 26

# Row 7
$ws.Range("H7").Value = 3.4
$ws.Range("I7").Value = 1.8
$ws.Range("J7").Value = 1.07
$ws.Range("K7").Value = 9
$ws.Range("L7").Value = 1.36
$ws.Range("M7").Value = 3
$ws.Range("N7").Value = 2.15
$ws.Range("O7").Value = 1.67
$ws.Range("R7").Value = 2
$ws.Range("S7").Value = 1.75
$ws.Range("Z7").Value = 8
$ws.Range("AE7").Value = 6
$ws.Range("AF7").Value = 7.5

# Row 8
$ws.Range("G8").Value = 2.3
$ws.Range("H8").Value = 3.3
$ws.Range("I8").Value = 3.1
$ws.Range("N8").Value = 2.25
$ws.Range("O8").Value = 1.62
$ws.Range("T8").Value = 7
$ws.Range("U8").Value = 10
$ws.Range("W8").Value = 21
$ws.Range("X8").Value = 21
$ws.Range("Z8").Value = 8
$ws.Range("AC8").Value = 51
$ws.Range("AD8").Value = 351
$ws.Range("AE8").Value = 8
$ws.Range("AF8").Value = 15
$ws.Range("AG8").Value = 11
$ws.Range("AH8").Value = 34
$ws.Range("AI8").Value = 26
